$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: G22 changes from nitem101 -> nitem101x, and H22:K22 are added
$ws.Range("G22").Value = "nitem101x"
$ws.Range("H22").Value = "nitem101x"
$ws.Range("I22").Value = "nitem102x"
$ws.Range("J22").Value = "nitem103x"
$ws.Range("K22").Value = "nitem104x"

# Row 23: G23 changes from nitem102 -> nitem102x, and H23 is added
$ws.Range("G23").Value = "nitem102x"
$ws.Range("H23").Value = "nitem102x"

# Row 24: G24 changes from nitem103 -> nitem103x, and H24 is added
$ws.Range("G24").Value = "nitem103x"
$ws.Range("H24").Value = "nitem103x"

# Row 25: G25 changes from nitem104 -> nitem104x, and H25 is added
$ws.Range("G25").Value = "nitem104x"
$ws.Range("H25").Value = "nitem104x"
